$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D..G to E..H)
$ws.Columns("D").EntireColumn.Insert()

# --- Header row ---
$ws.Range("C1").Value = "OPEX [EUR/kWh]"
$ws.Range("D1").Value = "OPEX [EUR/GWh]"

# --- Row 2 (BESS) ---
$ws.Range("B2").Value = 238
$ws.Range("C2").Formula = "=0.03*B2"
$ws.Range("D2").Formula = "=C2*1000000"
$ws.Range("D2").NumberFormat = "_-* #,##0.00_-;\-* #,##0.00_-;_-* ""-""??_-;_-@_-"
$ws.Range("F2").Style = "Percent"
$ws.Range("G2").Style = "Percent"
$ws.Range("G2").NumberFormat = "0.00%"
$ws.Range("H2").Formula = "=(B2*G2)*1000000"

# --- Row 3 (Pumped Hydro) ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 12.4
$ws.Range("D3").Formula = "=C3*1000000"
$ws.Range("D3").NumberFormat = "_-* #,##0.00_-;\-* #,##0.00_-;_-* ""-""??_-;_-@_-"
$ws.Range("E3:G3").ClearContents()
$ws.Range("H3").Value = 0

# --- Selection ---
$ws.Range("D6").Select()
